# budgetBills.xlsx: rename the "prise" (typo) header to "amount"
# and move the active selection from C8 to H2, matching the author's
# manual edit when he adjusted the test fixture for the new
# filterAccountBillsNotInBudgetBills() logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "amount"
$ws.Range("H2").Select()
